# Commit: "FIX when a named range is not in the excel ref."
# Adds a regression-test cell on Sheet1: D5 = INDIRECT("invalid"),
# which resolves to an invalid/undefined reference and therefore
# evaluates to the #REF! error. Mirrors the sibling INDIRECT() cells
# already present in column D (D2, D3, D4).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("D5").Formula = '=INDIRECT("invalid")'

# Typing a formula into D5 and pressing Enter leaves the active cell on
# the next row down (D6) - matches the saved selection state.
$ws.Range("D6").Select()
